$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new star/planet inputs (smaller stellar radius from 2023 paper) ---
# B2 (lat) stays the same value, but gets re-styled like the "Input" cells
# (bold font, orange fill, full border box) while keeping 3-decimal formatting.
$ws.Range("C2").Value = 0.59940000000000004

# --- Row 5: new simulation with different inclination (B5) and same stellar params as row 2 ---
$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 0.59940000000000004

# --- Row 10: manually pasted copy of the K5 result ---
$ws.Range("K10").Value = 87.808824582594895

# --- Re-style B2 to match the bordered/highlighted "Input" look (same as D2), ---
# --- and switch its number format used for B2 from 5 decimals to 3 decimals ---
$ws.Range("D2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").NumberFormat = "0.000"
$leftBorder = $ws.Range("B2").Borders.Item(7)
$leftBorder.LineStyle = 1
$leftBorder.ColorIndex = 1
$ws.Application.CutCopyMode = $false

# --- Selection moved ---
$ws.Range("H25").Select()
